$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing labels (shared string text changes from "Reel N" to "Image N")
$ws.Range("A6").Value = "Image 1"
$ws.Range("A7").Value = "Image 2"
$ws.Range("A8").Value = "Image 3"

# Update D8 value 383 -> 384
$ws.Range("D8").Value = 384

# Add new rows 9-12
$ws.Range("A9").Value = "Jackpot Text"
$ws.Range("B9").Value = 65
$ws.Range("C9").Value = 19
$ws.Range("D9").Value = 293
$ws.Range("E9").Value = 107

$ws.Range("A10").Value = "Credits Text"
$ws.Range("B10").Value = 65
$ws.Range("C10").Value = 19
$ws.Range("D10").Value = 194
$ws.Range("E10").Value = 303

$ws.Range("A11").Value = "Bet Text"
$ws.Range("B11").Value = 65
$ws.Range("C11").Value = 19
$ws.Range("D11").Value = 291
$ws.Range("E11").Value = 303

$ws.Range("A12").Value = "Result Text"
$ws.Range("B12").Value = 65
$ws.Range("C12").Value = 19
$ws.Range("D12").Value = 390
$ws.Range("E12").Value = 303

# copy style from A8/B8 row to new rows for consistent formatting
$ws.Range("A8").Copy()
$ws.Range("A9:A12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B8:E8").Copy()
$ws.Range("B9:E9").PasteSpecial(-4122)
$ws.Range("B9:E9").Copy()
$ws.Range("B10:E12").PasteSpecial(-4122)

# Update selection to G10 (as indicated by the diff)
$ws.Range("G10").Select()
